$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Week 6 entries (Nathan's hours wk 5 per the commit message, recorded under week label 6)
# Row 20
$ws.Range("A20").Value = 6
$ws.Range("B20").Value = 43704
$ws.Range("C20").Value = 0.39583333333333331
$ws.Range("D20").Value = 43704
$ws.Range("E20").Value = 0.40972222222222227
$ws.Range("F20").Value = "Formal team meeting"

# Row 22 is filled in before Row 21 so the new shared strings are appended to the
# sharedStrings table in the same order Excel originally recorded them.
# Row 22
$ws.Range("A22").Value = 6
$ws.Range("B22").Value = 43707
$ws.Range("C22").Value = 0.54166666666666663
$ws.Range("D22").Value = 43707
$ws.Range("E22").Value = 0.625
$ws.Range("F22").Value = "Creating user database + authentication "

# Row 21
$ws.Range("A21").Value = 6
$ws.Range("B21").Value = 43706
$ws.Range("C21").Value = 0.83333333333333337
$ws.Range("D21").Value = 43706
$ws.Range("E21").Value = 0.875
$ws.Range("F21").Value = "Began registration form "

$wb.Save()
